$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 459.35483
$ws.Range("I80").Value = 312.8
$ws.Range("J80").Value = 725.8182
$ws.Range("K80").Value = 938.4000000000001
$ws.Range("L80").Value = 2177.4546
$ws.Range("M80").Value = 59.59999999999991
$ws.Range("N80").Value = -4173.4546

$ws.Range("H83").Value = 459.35483
$ws.Range("I83").Value = 312.8
$ws.Range("J83").Value = 725.8182
$ws.Range("K83").Value = 2815.2
$ws.Range("L83").Value = 6532.3638
$ws.Range("M83").Value = 2176.8
$ws.Range("N83").Value = -16516.3638

$ws.Range("H94").Value = 2444.4443
$ws.Range("I94").Value = 2444.4443
$ws.Range("K94").Value = 2444.4443
$ws.Range("M94").Value = -1993.4443

$ws.Range("H106").Value = 8000
$ws.Range("I106").Value = 8000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 8000
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -7369
$ws.Range("N106").ClearContents()

$ws.Range("H138").Value = 3404709.5
$ws.Range("J138").Value = 4277300.5
$ws.Range("L138").Value = 12831901.5
$ws.Range("N138").Value = -12842181.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H97").Value = 3125785.5
$ws.Range("I97").Value = 3677103
$ws.Range("K97").Value = 3677103
$ws.Range("M97").Value = -3676607

$ws.Range("H110").Value = 1466
$ws.Range("I110").Value = 1366.3334
$ws.Range("K110").Value = 1366.3334
$ws.Range("M110").Value = 678.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 42527.625
$ws.Range("I86").Value = 57800
$ws.Range("J86").Value = 27255.25
$ws.Range("K86").Value = 57800
$ws.Range("L86").Value = 27255.25
$ws.Range("M86").Value = -56677
$ws.Range("N86").Value = -29501.25

$ws.Range("H89").Value = 42527.625
$ws.Range("I89").Value = 57800
$ws.Range("J89").Value = 27255.25
$ws.Range("K89").Value = 289000
$ws.Range("L89").Value = 136276.25
$ws.Range("M89").Value = -283384
$ws.Range("N89").Value = -147508.25

$ws.Range("H99").Value = 1149.6364
$ws.Range("I99").Value = 1197.5
$ws.Range("J99").Value = 1076
$ws.Range("K99").Value = 1197.5
$ws.Range("L99").Value = 1076
$ws.Range("M99").Value = 300.5
$ws.Range("N99").Value = -4072

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2533.3333
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2533.3333
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2533.3333
$ws.Range("N16").Value = -3107.3333
$ws.Range("M16").ClearContents()

$ws.Range("H113").Value = 2533.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2533.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2533.3333
$ws.Range("N113").Value = -6873.3333
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 417059.16
$ws.Range("I11").Value = 182
$ws.Range("K11").Value = 546
$ws.Range("M11").Value = -406

$ws.Range("H106").Value = 3183.4167
$ws.Range("J106").Value = 3183.4167
$ws.Range("L106").Value = 9550.250100000001
$ws.Range("N106").Value = -11442.2501

$ws.Range("H131").Value = 1082.4615
$ws.Range("J131").Value = 1199.2954
$ws.Range("L131").Value = 3597.8862
$ws.Range("N131").Value = -13677.8862

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 45230
$ws.Range("I70").Value = 76289.28999999999
$ws.Range("J70").Value = 5700
$ws.Range("K70").Value = 76289.28999999999
$ws.Range("L70").Value = 5700
$ws.Range("M70").Value = -76019.28999999999
$ws.Range("N70").Value = -6240

$ws.Range("H73").Value = 45230
$ws.Range("I73").Value = 76289.28999999999
$ws.Range("J73").Value = 5700
$ws.Range("K73").Value = 76289.28999999999
$ws.Range("L73").Value = 5700
$ws.Range("M73").Value = -75353.28999999999
$ws.Range("N73").Value = -7572

$ws.Range("H80").Value = 3931.1765
$ws.Range("I80").Value = 3290
$ws.Range("J80").Value = 4016.6667
$ws.Range("K80").Value = 3290
$ws.Range("L80").Value = 4016.6667
$ws.Range("M80").Value = -2292
$ws.Range("N80").Value = -6012.6667

$ws.Range("H83").Value = 3931.1765
$ws.Range("I83").Value = 3290
$ws.Range("J83").Value = 4016.6667
$ws.Range("K83").Value = 16450
$ws.Range("L83").Value = 20083.3335
$ws.Range("M83").Value = -11458
$ws.Range("N83").Value = -30067.3335

$ws.Range("H97").Value = 1079.5
$ws.Range("J97").Value = 400.5
$ws.Range("L97").Value = 400.5
$ws.Range("N97").Value = -1392.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2410.3914
$ws.Range("I16").Value = 830.4286
$ws.Range("J16").Value = 19000
$ws.Range("K16").Value = 830.4286
$ws.Range("L16").Value = 19000
$ws.Range("M16").Value = -660.4286
$ws.Range("N16").Value = -19340

$ws.Range("H68").Value = 1700.5
$ws.Range("I68").Value = 1401
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1401
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -652
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 1700.5
$ws.Range("I71").Value = 1401
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7005
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -3261
$ws.Range("N71").Value = -17488

$ws.Range("H82").Value = 3286.9092
$ws.Range("I82").Value = 1175.7273
$ws.Range("J82").Value = 5398.091
$ws.Range("K82").Value = 1175.7273
$ws.Range("L82").Value = 5398.091
$ws.Range("M82").Value = -814.7273
$ws.Range("N82").Value = -6120.091

$ws.Range("H85").Value = 3286.9092
$ws.Range("I85").Value = 1175.7273
$ws.Range("J85").Value = 5398.091
$ws.Range("K85").Value = 1175.7273
$ws.Range("L85").Value = 5398.091
$ws.Range("M85").Value = 72.27269999999999
$ws.Range("N85").Value = -7894.091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 92787.59
$ws.Range("I132").Value = 92295.37
$ws.Range("J132").Value = 93279.82000000001
$ws.Range("K132").Value = 276886.11
$ws.Range("L132").Value = 279839.46
$ws.Range("M132").Value = -274356.11
$ws.Range("N132").Value = -284899.46
